$d = $word.ActiveDocument

# The document has one section whose primary (default) and first-page
# headers/footers each contain a single inline picture (the Pearson logo
# in the footers, the BTEC logo in the headers). The commit renames the
# inline pictures:
#   footers: image1.png -> image2.png
#   headers: image2.jpg -> image1.jpg
# This is done through the InlineShape.Name property, which is the
# Word object-model member backing <wp:docPr name="...">.

for ($si = 1; $si -le $d.Sections.Count; $si++) {
    $sec = $d.Sections.Item($si)

    for ($hi = 1; $hi -le 2; $hi++) {
        $hdr = $sec.Headers.Item($hi)
        if ($hdr.Exists) {
            for ($k = 1; $k -le $hdr.Range.InlineShapes.Count; $k++) {
                $shp = $hdr.Range.InlineShapes.Item($k)
                if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                    $shp.Name = "image1.jpg"
                }
            }
        }
    }

    for ($fi = 1; $fi -le 2; $fi++) {
        $ftr = $sec.Footers.Item($fi)
        if ($ftr.Exists) {
            for ($k = 1; $k -le $ftr.Range.InlineShapes.Count; $k++) {
                $shp = $ftr.Range.InlineShapes.Item($k)
                if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                    $shp.Name = "image2.png"
                }
            }
        }
    }
}
